$d = $word.ActiveDocument

# Paragraph layout before edit:
#   1: To Do List:
#   2: Implement mortgage into the turn
#   3: Make it so that the get out of jail free card can be reshuffled into the deck after use.
#   4: Bug Hunting
#   5: Even when you escape by rolling doubles, it repeats the jail dialogue until you pay the money
#   6: Luxury tax take money for every turn, including other players turns
#   7: (empty, holds the _GoBack bookmark)

# Remove paragraphs 2 and 3 ("Implement mortgage into the turn" and
# "Make it so that the get out of jail free card can be reshuffled into the
# deck after use.") in one go, which also merges "Bug Hunting" up into
# what was paragraph 2's slot.
$pStart = $d.Paragraphs.Item(2)
$pEnd = $d.Paragraphs.Item(3)
$d.Range($pStart.Range.Start, $pEnd.Range.End).Delete()

# Remove the paragraph that now follows "Bug Hunting":
# "Even when you escape by rolling doubles, it repeats the jail dialogue
# until you pay the money"
$pEscape = $d.Paragraphs.Item(3)
$pEscape.Range.Delete()

# Move the _GoBack bookmark from the trailing empty paragraph to the end of
# the "Bug Hunting" text (i.e. immediately before its paragraph mark).
#
# A collapsed Range placed exactly one position before a paragraph mark
# confuses this host's Bookmarks.Add, so insert a throwaway placeholder
# character there first (pushing the paragraph mark out by one), anchor the
# bookmark just before the placeholder, and then delete the placeholder.
$pBugHunting = $d.Paragraphs.Item(2)
$insertPoint = $d.Range($pBugHunting.Range.End - 1, $pBugHunting.Range.End - 1)
$insertPoint.InsertAfter("X")

$pBugHunting2 = $d.Paragraphs.Item(2)
$bmPos = $pBugHunting2.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($bmPos, $bmPos + 1).Delete()
